# Auto-generated Excel COM-interop edit script
# Applies the numeric corrections described in the commit diff
# for workbook "Asura_Profits" (sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 9842.75
$ws.Range("I92").Value = 9842.75
$ws.Range("K92").Value = 9842.75
$ws.Range("M92").Value = -8594.75
$ws.Range("H100").Value = 3463.5
$ws.Range("I100").Value = 3501.4285
$ws.Range("J100").Value = 3434
$ws.Range("K100").Value = 3501.4285
$ws.Range("L100").Value = 3434
$ws.Range("M100").Value = -2960.4285
$ws.Range("N100").Value = -4516
$ws.Range("H127").Value = 1165.0883
$ws.Range("I127").Value = 990
$ws.Range("J127").Value = 1167.7015
$ws.Range("K127").Value = 2970
$ws.Range("L127").Value = 3503.104499999999
$ws.Range("M127").Value = 1990
$ws.Range("N127").Value = -13423.1045
$ws.Range("H137").Value = 1408.4546
$ws.Range("I137").Value = 1298.8667
$ws.Range("J137").Value = 1643.2858
$ws.Range("K137").Value = 3896.6001
$ws.Range("L137").Value = 4929.857400000001
$ws.Range("M137").Value = -1346.6001
$ws.Range("N137").Value = -10029.8574
$ws.Range("H138").Value = 2969.2585
$ws.Range("I138").Value = 1852.4062
$ws.Range("J138").Value = 4343.846
$ws.Range("K138").Value = 5557.2186
$ws.Range("L138").Value = 13031.538
$ws.Range("M138").Value = -417.2186000000002
$ws.Range("N138").Value = -23311.538

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 68909.2
$ws.Range("I2").Value = 1261.25
$ws.Range("J2").Value = 93508.45
$ws.Range("K2").Value = 1261.25
$ws.Range("L2").Value = 93508.45
$ws.Range("M2").Value = -1148.25
$ws.Range("N2").Value = -93734.45
$ws.Range("H32").Value = 10011.464
$ws.Range("I32").Value = 10306.967
$ws.Range("K32").Value = 10306.967
$ws.Range("M32").Value = -10019.967
$ws.Range("H45").Value = 961.381
$ws.Range("I45").Value = 926.26666
$ws.Range("J45").Value = 1049.1666
$ws.Range("K45").Value = 926.26666
$ws.Range("L45").Value = 1049.1666
$ws.Range("M45").Value = -549.26666
$ws.Range("N45").Value = -1803.1666
$ws.Range("H61").Value = 2977.1035
$ws.Range("I61").Value = 2174.818
$ws.Range("J61").Value = 5498.5713
$ws.Range("K61").Value = 2174.818
$ws.Range("L61").Value = 5498.5713
$ws.Range("M61").Value = -1962.818
$ws.Range("N61").Value = -5922.5713
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H74").Value = 977.02856
$ws.Range("I74").Value = 845.1539
$ws.Range("J74").Value = 1358
$ws.Range("K74").Value = 845.1539
$ws.Range("L74").Value = 1358
$ws.Range("M74").Value = 28.84609999999998
$ws.Range("N74").Value = -3106
$ws.Range("H77").Value = 977.02856
$ws.Range("I77").Value = 845.1539
$ws.Range("J77").Value = 1358
$ws.Range("K77").Value = 4225.7695
$ws.Range("L77").Value = 6790
$ws.Range("M77").Value = 142.2304999999997
$ws.Range("N77").Value = -15526
$ws.Range("H102").Value = 335300
$ws.Range("I102").Value = 2900
$ws.Range("J102").Value = 501500
$ws.Range("K102").Value = 2900
$ws.Range("L102").Value = 501500
$ws.Range("M102").Value = -1278
$ws.Range("N102").Value = -504744
$ws.Range("H116").Value = 68909.2
$ws.Range("I116").Value = 1261.25
$ws.Range("J116").Value = 93508.45
$ws.Range("K116").Value = 1261.25
$ws.Range("L116").Value = 93508.45
$ws.Range("M116").Value = 1032.75
$ws.Range("N116").Value = -98096.45
$ws.Range("H132").Value = 4387.281
$ws.Range("I132").Value = 6256.3076
$ws.Range("J132").Value = 2819.7097
$ws.Range("K132").Value = 18768.9228
$ws.Range("L132").Value = 8459.1291
$ws.Range("M132").Value = -16238.9228
$ws.Range("N132").Value = -13519.1291
$ws.Range("H136").Value = 2977.1035
$ws.Range("I136").Value = 2174.818
$ws.Range("J136").Value = 5498.5713
$ws.Range("K136").Value = 6524.454000000001
$ws.Range("L136").Value = 16495.7139
$ws.Range("M136").Value = -3974.454000000001
$ws.Range("N136").Value = -21595.7139

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 68909.2
$ws.Range("I3").Value = 1261.25
$ws.Range("J3").Value = 93508.45
$ws.Range("K3").Value = 1261.25
$ws.Range("L3").Value = 93508.45
$ws.Range("M3").Value = -1147.25
$ws.Range("N3").Value = -93736.45
$ws.Range("H86").Value = 102725.8
$ws.Range("I86").Value = 3157.7856
$ws.Range("J86").Value = 335051.16
$ws.Range("K86").Value = 3157.7856
$ws.Range("L86").Value = 335051.16
$ws.Range("M86").Value = -2034.7856
$ws.Range("N86").Value = -337297.16
$ws.Range("H89").Value = 102725.8
$ws.Range("I89").Value = 3157.7856
$ws.Range("J89").Value = 335051.16
$ws.Range("K89").Value = 15788.928
$ws.Range("L89").Value = 1675255.8
$ws.Range("M89").Value = -10172.928
$ws.Range("N89").Value = -1686487.8
$ws.Range("H108").Value = 54684
$ws.Range("J108").Value = 54684
$ws.Range("L108").Value = 54684
$ws.Range("N108").Value = -62364
$ws.Range("H122").Value = 55740.11
$ws.Range("J122").Value = 55740.11
$ws.Range("L122").Value = 55740.11
$ws.Range("N122").Value = -65540.11
$ws.Range("H134").Value = 6272.638
$ws.Range("I134").Value = 2094.8333
$ws.Range("J134").Value = 19945.455
$ws.Range("K134").Value = 6284.499899999999
$ws.Range("L134").Value = 59836.36500000001
$ws.Range("M134").Value = -3749.499899999999
$ws.Range("N134").Value = -64906.36500000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 27670
$ws.Range("I14").Value = 27505
$ws.Range("J14").Value = 28000
$ws.Range("K14").Value = 27505
$ws.Range("L14").Value = 28000
$ws.Range("M14").Value = -27335
$ws.Range("N14").Value = -28340
$ws.Range("H31").Value = 2706.6924
$ws.Range("I31").Value = 1798.375
$ws.Range("J31").Value = 4160
$ws.Range("K31").Value = 1798.375
$ws.Range("L31").Value = 4160
$ws.Range("M31").Value = -1503.375
$ws.Range("N31").Value = -4750
$ws.Range("H34").Value = 2706.6924
$ws.Range("I34").Value = 1798.375
$ws.Range("J34").Value = 4160
$ws.Range("K34").Value = 1798.375
$ws.Range("L34").Value = 4160
$ws.Range("M34").Value = -1596.375
$ws.Range("N34").Value = -4564
$ws.Range("H58").Value = 741821.8
$ws.Range("I58").Value = 975672.3
$ws.Range("J58").Value = 1295.0834
$ws.Range("K58").Value = 975672.3
$ws.Range("L58").Value = 1295.0834
$ws.Range("M58").Value = -975469.3
$ws.Range("N58").Value = -1701.0834
$ws.Range("H132").Value = 266335.97
$ws.Range("I132").Value = 356647.25
$ws.Range("J132").Value = 2349.077
$ws.Range("K132").Value = 1069941.75
$ws.Range("L132").Value = 7047.231000000001
$ws.Range("M132").Value = -1067411.75
$ws.Range("N132").Value = -12107.231
$ws.Range("H134").Value = 1536.1025
$ws.Range("I134").Value = 1225.7037
$ws.Range("J134").Value = 2234.5
$ws.Range("K134").Value = 3677.1111
$ws.Range("L134").Value = 6703.5
$ws.Range("M134").Value = -1142.1111
$ws.Range("N134").Value = -11773.5
$ws.Range("H136").Value = 741821.8
$ws.Range("I136").Value = 975672.3
$ws.Range("J136").Value = 1295.0834
$ws.Range("K136").Value = 2927016.9
$ws.Range("L136").Value = 3885.2502
$ws.Range("M136").Value = -2924466.9
$ws.Range("N136").Value = -8985.2502

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 152.63158
$ws.Range("I38").Value = 131.33333
$ws.Range("J38").Value = 232.5
$ws.Range("K38").Value = 393.99999
$ws.Range("L38").Value = 697.5
$ws.Range("M38").Value = -46.99998999999997
$ws.Range("N38").Value = -1391.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3582.1904
$ws.Range("I132").Value = 3413.3235
$ws.Range("J132").Value = 4299.875
$ws.Range("K132").Value = 10239.9705
$ws.Range("L132").Value = 12899.625
$ws.Range("M132").Value = -7709.970499999999
$ws.Range("N132").Value = -17959.625

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4047.0588
$ws.Range("I40").Value = 4144.5557
$ws.Range("J40").Value = 3937.375
$ws.Range("K40").Value = 4144.5557
$ws.Range("L40").Value = 3937.375
$ws.Range("M40").Value = -4008.5557
$ws.Range("N40").Value = -4209.375
$ws.Range("H55").Value = 603.1539
$ws.Range("I55").Value = 296.25
$ws.Range("J55").Value = 739.55554
$ws.Range("K55").Value = 296.25
$ws.Range("L55").Value = 739.55554
$ws.Range("M55").Value = -123.25
$ws.Range("N55").Value = -1085.55554
$ws.Range("H61").Value = 12381.714
$ws.Range("I61").Value = 13348.211
$ws.Range("J61").Value = 3200
$ws.Range("K61").Value = 13348.211
$ws.Range("L61").Value = 3200
$ws.Range("M61").Value = -13146.211
$ws.Range("N61").Value = -3604
$ws.Range("H113").Value = 12381.714
$ws.Range("I113").Value = 13348.211
$ws.Range("J113").Value = 3200
$ws.Range("K113").Value = 13348.211
$ws.Range("L113").Value = 3200
$ws.Range("M113").Value = -11178.211
$ws.Range("N113").Value = -7540
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H136").Value = 3467.923
$ws.Range("I136").Value = 2564.5
$ws.Range("J136").Value = 5500.625
$ws.Range("K136").Value = 7693.5
$ws.Range("L136").Value = 16501.875
$ws.Range("M136").Value = -5143.5
$ws.Range("N136").Value = -21601.875

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 50000
$ws.Range("J76").Value = 50000
$ws.Range("L76").Value = 50000
$ws.Range("N76").Value = -50630
$ws.Range("H79").Value = 50000
$ws.Range("J79").Value = 50000
$ws.Range("L79").Value = 50000
$ws.Range("N79").Value = -52184
$ws.Range("H107").Value = 728
$ws.Range("I107").Value = 716.6667
$ws.Range("J107").Value = 745
$ws.Range("K107").Value = 2150.0001
$ws.Range("L107").Value = 2235
$ws.Range("M107").Value = -230.0001000000002
$ws.Range("N107").Value = -6075
$ws.Range("H126").Value = 8059.778
$ws.Range("I126").Value = 9812.929
$ws.Range("J126").Value = 1923.75
$ws.Range("K126").Value = 29438.787
$ws.Range("L126").Value = 5771.25
$ws.Range("M126").Value = -26968.787
$ws.Range("N126").Value = -10711.25
$ws.Range("H132").Value = 1023.0345
$ws.Range("I132").Value = 813.1875
$ws.Range("J132").Value = 2030.3
$ws.Range("K132").Value = 2439.5625
$ws.Range("L132").Value = 6090.9
$ws.Range("M132").Value = 90.4375
$ws.Range("N132").Value = -11150.9
$ws.Range("H133").Value = 29800
$ws.Range("J133").Value = 29800
$ws.Range("L133").Value = 29800
$ws.Range("N133").Value = -39920
$ws.Range("H136").Value = 1173
$ws.Range("I136").Value = 1396.0385
$ws.Range("K136").Value = 4188.1155
$ws.Range("M136").Value = -1638.1155
